$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3085.6775
$ws.Range("I141").Value = 1393.0476
$ws.Range("K141").Value = 4179.142800000001
$ws.Range("M141").Value = 1000.857199999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1026.16
$ws.Range("I32").Value = 794.13336
$ws.Range("J32").Value = 3114.4
$ws.Range("K32").Value = 794.13336
$ws.Range("L32").Value = 3114.4
$ws.Range("M32").Value = -507.13336
$ws.Range("N32").Value = -3688.4
$ws.Range("H37").Value = 11601.637
$ws.Range("I37").Value = 5447.5
$ws.Range("J37").Value = 28012.666
$ws.Range("K37").Value = 5447.5
$ws.Range("L37").Value = 28012.666
$ws.Range("M37").Value = -5174.5
$ws.Range("N37").Value = -28558.666
$ws.Range("H61").Value = 3244295
$ws.Range("I61").Value = 2293623
$ws.Range("J61").Value = 8405087
$ws.Range("K61").Value = 2293623
$ws.Range("L61").Value = 8405087
$ws.Range("M61").Value = -2293411
$ws.Range("N61").Value = -8405511
$ws.Range("H136").Value = 3244295
$ws.Range("I136").Value = 2293623
$ws.Range("J136").Value = 8405087
$ws.Range("K136").Value = 6880869
$ws.Range("L136").Value = 25215261
$ws.Range("M136").Value = -6878319
$ws.Range("N136").Value = -25220361

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524
$ws.Range("H134").Value = 11765424
$ws.Range("I134").Value = 11628612
$ws.Range("K134").Value = 34885836
$ws.Range("M134").Value = -34883301

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 61522
$ws.Range("J106").Value = 61522
$ws.Range("L106").Value = 61522
$ws.Range("N106").Value = -64046
$ws.Range("H108").Value = 28314
$ws.Range("J108").Value = 28314
$ws.Range("L108").Value = 28314
$ws.Range("N108").Value = -35994

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3156.2
$ws.Range("I94").Value = 350
$ws.Range("J94").Value = 5027
$ws.Range("K94").Value = 1050
$ws.Range("L94").Value = 15081
$ws.Range("M94").Value = -374
$ws.Range("N94").Value = -16433
$ws.Range("H95").Value = 13000
$ws.Range("J95").Value = 13000
$ws.Range("L95").Value = 39000
$ws.Range("N95").Value = -43118
$ws.Range("H96").Value = 2375
$ws.Range("I96").Value = 3412.5
$ws.Range("J96").Value = 2144.4443
$ws.Range("K96").Value = 10237.5
$ws.Range("L96").Value = 6433.3329
$ws.Range("M96").Value = -8178.5
$ws.Range("N96").Value = -10551.3329
$ws.Range("H97").Value = 726.2222
$ws.Range("I97").Value = 800
$ws.Range("J97").Value = 667.2
$ws.Range("K97").Value = 2400
$ws.Range("L97").Value = 2001.6
$ws.Range("M97").Value = -1904
$ws.Range("N97").Value = -2993.6
$ws.Range("H98").Value = 778
$ws.Range("I98").Value = 778
$ws.Range("K98").Value = 2334
$ws.Range("M98").Value = -836
$ws.Range("H99").Value = 16011.25
$ws.Range("I99").Value = 17015
$ws.Range("K99").Value = 51045
$ws.Range("M99").Value = -48799
$ws.Range("H100").Value = 4973.6665
$ws.Range("J100").Value = 4973.6665
$ws.Range("L100").Value = 14920.9995
$ws.Range("N100").Value = -16542.9995
$ws.Range("H101").Value = 9999.429
$ws.Range("J101").Value = 9999.429
$ws.Range("L101").Value = 29998.287
$ws.Range("N101").Value = -34866.287
$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 3000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 9000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -6566
$ws.Range("N102").ClearContents()
$ws.Range("H103").Value = 1533.2858
$ws.Range("J103").Value = 1633.25
$ws.Range("L103").Value = 4899.75
$ws.Range("N103").Value = -6657.75
$ws.Range("H104").Value = 3000
$ws.Range("J104").Value = 3000
$ws.Range("L104").Value = 9000
$ws.Range("N104").Value = -14242
$ws.Range("H105").Value = 5089.125
$ws.Range("J105").Value = 4960
$ws.Range("L105").Value = 14880
$ws.Range("N105").Value = -20122
$ws.Range("H106").Value = 5994.375
$ws.Range("J106").Value = 5994.375
$ws.Range("L106").Value = 17983.125
$ws.Range("N106").Value = -19875.125
$ws.Range("H107").Value = 435279.5
$ws.Range("I107").Value = 854908.25
$ws.Range("J107").Value = 1180.7931
$ws.Range("K107").Value = 2564724.75
$ws.Range("L107").Value = 3542.379300000001
$ws.Range("M107").Value = -2562804.75
$ws.Range("N107").Value = -7382.379300000001
$ws.Range("H108").Value = 1228.4
$ws.Range("I108").Value = 469.14285
$ws.Range("K108").Value = 1407.42855
$ws.Range("M108").Value = 1472.57145
$ws.Range("H109").Value = 3444.7805
$ws.Range("I109").Value = 1159.5714
$ws.Range("J109").Value = 3915.2646
$ws.Range("K109").Value = 3478.7142
$ws.Range("L109").Value = 11745.7938
$ws.Range("M109").Value = -2438.7142
$ws.Range("N109").Value = -13825.7938
$ws.Range("H120").Value = 111111950
$ws.Range("I120").Value = 111111950
$ws.Range("K120").Value = 333335850
$ws.Range("M120").Value = -333331012

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8282570.5
$ws.Range("I132").Value = 4792234
$ws.Range("J132").Value = 29224590
$ws.Range("K132").Value = 14376702
$ws.Range("L132").Value = 87673770
$ws.Range("M132").Value = -14374172
$ws.Range("N132").Value = -87678830

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 114058.89
$ws.Range("I16").Value = 145433
$ws.Range("J16").Value = 4249.5
$ws.Range("K16").Value = 145433
$ws.Range("L16").Value = 4249.5
$ws.Range("M16").Value = -145263
$ws.Range("N16").Value = -4589.5
$ws.Range("H40").Value = 3587974.8
$ws.Range("I40").Value = 5294124.5
$ws.Range("J40").Value = 5060.3
$ws.Range("K40").Value = 5294124.5
$ws.Range("L40").Value = 5060.3
$ws.Range("M40").Value = -5293988.5
$ws.Range("N40").Value = -5332.3
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26622
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -83112

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 916066.7
$ws.Range("I132").Value = 353280.97
$ws.Range("J132").Value = 2757910.8
$ws.Range("K132").Value = 1059842.91
$ws.Range("L132").Value = 8273732.399999999
$ws.Range("M132").Value = -1057312.91
$ws.Range("N132").Value = -8278792.399999999
